$p = $ppt.ActivePresentation

# Slide 2: "den 13.05.2017, 24 Uhr mit, wie Sie Ihre Gruppen" -> 2018
# This run sits between two manual line breaks (<a:br/>), so it is
# addressed via TextRange.Characters(start, length) -- found dynamically
# by locating the old substring inside the paragraph's plain text -- rather
# than Paragraphs().Runs(), which misbehaves for runs flanked by line breaks
# on both sides.
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(5)
$tr2 = $shp2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(1)
$oldRun2 = "den 13.05.2017, 24 Uhr mit, wie Sie Ihre Gruppen"
$newRun2 = "den 13.05.2018, 24 Uhr mit, wie Sie Ihre Gruppen"
$idx2 = $para2.Text.IndexOf($oldRun2)
$tr2.Characters($para2.Start + $idx2, $oldRun2.Length).Text = $newRun2

# Slide 3: hausarbeit-w15c- -> hausarbeit-w16c-
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(5)
$tr3 = $shp3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(2)
$para3.Runs(2).Text = "https://github.com/nordakademie-einfuehrung-java/hausarbeit-w16c-"

# Slide 4: two date/text changes
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(5)
$tr4 = $shp4.TextFrame.TextRange
$tr4.Paragraphs(1).Runs(1).Text = "Das Thema zur Hausarbeit wird Ihnen spätestens am 17.05.2018 um 0:00 Uhr per E-Mail an den Zenturien-Verteiler zugestellt"
$tr4.Paragraphs(4).Runs(1).Text = "Während der Vorlesung vom 24.05.2018 können Fragen zur Aufgabenstellung geklärt werden, die nicht bereits per E-Mail behandelt wurden"

# Slide 5: replace last bullet text and insert two new paragraphs before it
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(5)
$tr5 = $shp5.TextFrame.TextRange
$tr5.Paragraphs(1).Runs(1).Text = "Die Hausarbeit ist am 13.06.2018 um 23:59 Uhr per E-Mail abzugeben"

# Insert two new paragraphs after paragraph 1 (current first bullet)
$para1 = $tr5.Paragraphs(1)
$para1.InsertAfter("`r" + "`r" + "In der letzten Vorlesung am 14.06.2018 werden die Ergebnisse von jeder Gruppe kurz vorgestellt") | Out-Null
